$wb = $excel.ActiveWorkbook

function Add-TradeRow68($ws) {
    $ws.Cells.Item(68, 1).Value = 67
    $ws.Cells.Item(68, 2).Value = "'2026-02-17"
    $ws.Cells.Item(68, 3).Value = "08:49:12"
    $ws.Cells.Item(68, 4).Value = "MarketMaking"
    $ws.Cells.Item(68, 5).Value = "DOWN"
    $ws.Cells.Item(68, 6).Value = 0.01
    $ws.Cells.Item(68, 7).Value = ""
    $ws.Cells.Item(68, 8).Value = "OPEN"
    $ws.Cells.Item(68, 9).Value = 0
    $ws.Cells.Item(68, 10).Value = 0
    $ws.Cells.Item(68, 11).Value = 99.67660198355652
    $ws.Cells.Item(68, 12).Value = 0
    $ws.Cells.Item(68, 13).Value = 0
    $ws.Cells.Item(68, 14).Value = 0.6
    $ws.Cells.Item(68, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(68, 16).Value = ""
    $ws.Cells.Item(68, 17).Value = 0
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow68 $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow68 $wsMarketMaking
